# Update gh-pages to output generated at 456a3b4
#
# Applies the same edit to both the "展览" sheet and the "全部类型" sheet
# (they carry identical data in this workbook):
#   1. Insert a new row at position 18 for a newly-listed event
#      (南昌·次元星球动漫游戏展), pushing the previous rows 18-39 down to 19-40.
#   2. Refresh the "想去人数" (F column) interest counters that ticked up for
#      a number of existing events.

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet($ws) {
    # --- 1. Insert the new row for 南昌·次元星球动漫游戏展 at row 18 -------------
    $ws.Rows("18:18").Insert()

    # Column A carries a bordered/bold/centered style (same as every other
    # row's index cell) - clone it from the row above, then overwrite the
    # value with the correct running index.
    $ws.Range("A17").Copy($ws.Range("A18"))
    $ws.Range("A18").Value = 17

    # Column B holds a plain text date string ("2024-07-06"). Force text
    # formatting first so Excel doesn't coerce it into a date serial, then
    # drop back to the Normal style so no stray number format sticks to the
    # cell (matching the unstyled cells used throughout the sheet).
    $ws.Range("B18").NumberFormat = "@"
    $ws.Range("B18").Value = "2024-07-06"
    $ws.Range("B18").Style = "Normal"

    $ws.Range("C18").Value = "南昌·次元星球动漫游戏展"
    $ws.Range("D18").Value = "龙蟠街666号融创茂1层 融创茂"
    $ws.Range("E18").Value = "2024.07.06 10:00-07.06 17:00"
    $ws.Range("F18").Value = 1
    $ws.Range("G18").Value = 50
    $ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86405"
    $ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg"

    # --- 2. Refresh "想去人数" (column F) counts ------------------------------
    # Rows 1-17 are untouched by the insert above, so these keep their
    # original row numbers.
    $ws.Range("F3").Value = 3119
    $ws.Range("F5").Value = 120
    $ws.Range("F7").Value = 1664
    $ws.Range("F8").Value = 1620
    $ws.Range("F9").Value = 56

    # Rows that used to be 19-36 are now 20-37 after the insert.
    $ws.Range("F20").Value = 15
    $ws.Range("F21").Value = 45
    $ws.Range("F24").Value = 181
    $ws.Range("F26").Value = 21
    $ws.Range("F28").Value = 73
    $ws.Range("F29").Value = 2094
    $ws.Range("F30").Value = 6
    $ws.Range("F33").Value = 199
    $ws.Range("F37").Value = 338
}

# Sheet "展览" (worksheet 1) and sheet "全部类型" (worksheet 4) both hold the
# same table and both need the identical update.
Update-ExpoSheet $wb.Worksheets.Item(1)
Update-ExpoSheet $wb.Worksheets.Item(4)
